$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume data scraped this run.
# D-column price strings are set with a leading apostrophe so Excel
# keeps them as literal text (matching the original inlineStr cells)
# instead of silently reinterpreting them as numbers; Style is reset
# back to Normal afterwards so no stray number-format style sticks.

$c = $ws.Range("D2")
$c.Value = "'26.399.46"
$c.Style = "Normal"
$ws.Range("E2").Value = "  -0.41%  "

$c = $ws.Range("D3")
$c.Value = "'1.722.54"
$c.Style = "Normal"
$ws.Range("E3").Value = "  -0.68%  "

$c = $ws.Range("D4")
$c.Value = "'0.9977"
$c.Style = "Normal"
$ws.Range("E4").Value = "  -0.17%  "

$c = $ws.Range("D5")
$c.Value = "'242.68"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -1.74%  "

$c = $ws.Range("D6")
$c.Value = "'0.9988"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -0.10%  "

$ws.Range("E7").Value = "  +0.30%  "

$c = $ws.Range("D8")
$c.Value = "'0.2605"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -2.24%  "

$c = $ws.Range("D9")
$c.Value = "'0.06187"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -0.52%  "

$c = $ws.Range("D10")
$c.Value = "'1.723.27"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -0.83%  "

$c = $ws.Range("D11")
$c.Value = "'0.06982"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -1.02%  "

$c = $ws.Range("D12")
$c.Value = "'15.54"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -1.43%  "

$c = $ws.Range("D13")
$c.Value = "'4.517"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -1.92%  "

$c = $ws.Range("D14")
$c.Value = "'0.5990"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -1.84%  "

$c = $ws.Range("D15")
$c.Value = "'77.07"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -0.43%  "

$c = $ws.Range("D16")
$c.Value = "'0.9986"
$c.Style = "Normal"

$c = $ws.Range("D17")
$c.Value = "'26.397.29"
$c.Style = "Normal"

$c = $ws.Range("D18")
$c.Value = "'0.9978"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -0.17%  "

$c = $ws.Range("D19")
$c.Value = "'0.000007135"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -1.17%  "

$ws.Range("E20").Value = "  -2.11%  "

$c = $ws.Range("D21")
$c.Value = "'1.944.78"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -0.71%  "

$c = $ws.Range("D22")
$c.Value = "'4.443"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -2.09%  "

$c = $ws.Range("D23")
$c.Value = "'8.501"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -2.98%  "

$c = $ws.Range("D24")
$c.Value = "'5.092"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -2.90%  "

$c = $ws.Range("D25")
$c.Value = "'137.91"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -1.11%  "

$c = $ws.Range("D26")
$c.Value = "'15.24"
$c.Style = "Normal"

$c = $ws.Range("D27")
$c.Value = "'1.403"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -0.63%  "

$c = $ws.Range("D28")
$c.Value = "'106.32"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -1.68%  "

$c = $ws.Range("D29")
$c.Value = "'1.735"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -2.66%  "

$ws.Range("E30").Value = "  -2.05%  "

$c = $ws.Range("D31")
$c.Value = "'0.08033"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -0.35%  "

$c = $ws.Range("D32")
$c.Value = "'3.651"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -1.13%  "

$c = $ws.Range("D33")
$c.Value = "'0.04485"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -2.06%  "

$c = $ws.Range("D34")
$c.Value = "'2.604"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -0.38%  "

$c = $ws.Range("D35")
$c.Value = "'0.9970"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -1.58%  "

$c = $ws.Range("D36")
$c.Value = "'0.6228"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -2.45%  "

$c = $ws.Range("D37")
$c.Value = "'0.9227"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +1.63%  "

$c = $ws.Range("D38")
$c.Value = "'1.963"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -3.75%  "

$c = $ws.Range("D39")
$c.Value = "'2.387"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -0.58%  "

$ws.Range("E40").Value = "  -0.58%  "

$c = $ws.Range("D41")
$c.Value = "'0.01478"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -2.05%  "

$c = $ws.Range("D42")
$c.Value = "'99.91"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -1.03%  "

$c = $ws.Range("D43")
$c.Value = "'5.409"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -0.31%  "

$c = $ws.Range("D44")
$c.Value = "'0.3841"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -1.59%  "

$c = $ws.Range("D45")
$c.Value = "'6.917"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -0.82%  "

$ws.Range("E46").Value = "  -1.92%  "

$c = $ws.Range("D47")
$c.Value = "'0.05367"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -0.46%  "

$c = $ws.Range("D48")
$c.Value = "'30.42"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -0.59%  "

$c = $ws.Range("D49")
$c.Value = "'7.696"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -1.31%  "

$c = $ws.Range("D50")
$c.Value = "'51.16"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -0.55%  "

$c = $ws.Range("D51")
$c.Value = "'1.218"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -2.73%  "

